$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -2
$ws.Range("F3").Value = -3
$ws.Range("F5").Value = -1
$ws.Range("F8").Value = 0
$ws.Range("F9").Value = -2
$ws.Range("F10").Value = -2
$ws.Range("F12").Value = -5
$ws.Range("F13").Value = 4
$ws.Range("F14").Value = 5
$ws.Range("F15").Value = 8
$ws.Range("F16").Value = 2
$ws.Range("F17").Value = -6
$ws.Range("F18").Value = 1
$ws.Range("F21").Value = 0
$ws.Range("F22").Value = -4
$ws.Range("F23").Value = -3
$ws.Range("F24").Value = 3
$ws.Range("F25").Value = 3
